$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 471, shifting all existing rows (471-553) down to (472-554)
$ws.Rows.Item(471).Insert()

# Populate the newly inserted row 471 with the new record
$ws.Cells.Item(471, 1).Value = 9
$ws.Cells.Item(471, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(471, 3).Value = "Metropolitana"
$ws.Cells.Item(471, 4).Value = 44637
$ws.Cells.Item(471, 5).Value = 13
$ws.Cells.Item(471, 6).Value = 100112040
$ws.Cells.Item(471, 7).Value = "Cilantro"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 220
$ws.Cells.Item(471, 11).Value = 12000
$ws.Cells.Item(471, 12).Value = 14000
$ws.Cells.Item(471, 13).Value = 13091
$ws.Cells.Item(471, 14).Value = "$/docena de atados"
$ws.Cells.Item(471, 15).Value = "Región Metropolitana"
$ws.Cells.Item(471, 16).Value = 4364
$ws.Cells.Item(471, 17).Value = 3
$ws.Cells.Item(471, 18).Value = "Hortaliza"

# Match the date-formatted style used by other rows in column D
$ws.Cells.Item(471, 4).NumberFormat = $ws.Cells.Item(472, 4).NumberFormat
